# This script re-applies a corrected set of match rows in the "Wales Premier
# League" sheet. Several rows that were previously mismatched (each row's
# id in column A matched to the wrong match's betting data in columns
# B..AC) are rotated into their correct positions. Column A (the row's
# sequential id) and the row position itself stay fixed; columns B through
# AC (match id, teams, scores, odds, etc.) move between rows according to
# the mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row. For every destination row, the
# values that should end up there (in columns B..AC) are the values that
# currently live (before this edit) in the source row.
$rowMap = @{
    22 = 23; 23 = 24; 24 = 22;
    26 = 27; 27 = 28; 28 = 26;
    29 = 31; 30 = 29; 31 = 30;
    53 = 54; 54 = 53;
    64 = 65; 65 = 66; 66 = 64;
    87 = 89; 88 = 87; 89 = 88;
    92 = 94; 93 = 92; 94 = 95; 95 = 93;
    98 = 99; 99 = 100; 100 = 98;
    104 = 105; 105 = 104;
    108 = 109; 109 = 108;
    119 = 120; 120 = 121; 121 = 119;
    129 = 132; 130 = 131; 131 = 133; 132 = 130; 133 = 129;
    140 = 141; 141 = 140;
}

$firstCol = 2   # column B
$lastCol  = 29  # column AC

# 1) Snapshot the current ("before") values of columns B..AC for every row
#    that takes part in the rotation, so that writes to one row never
#    clobber data still needed for another row later in the loop.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowData = @()
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowData += ,$ws.Cells.Item($r, $col).Value()
    }
    $snapshot[$r] = $rowData
}

# 2) Write each destination row's columns B..AC from the snapshot of its
#    mapped source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $data = $snapshot[$srcRow]
    for ($i = 0; $i -lt $data.Count; $i++) {
        $col = $firstCol + $i
        $ws.Cells.Item($destRow, $col).Value = $data[$i]
    }
}
